$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - Min 6 purchases: fill in the previously-missing raw counts
$ws.Range("D21").Value = 1063
$ws.Range("E21").Value = 2795
$ws.Range("F21").Value = 1085
$ws.Range("G21").Value = 2773
$ws.Range("H21").Value = 936
$ws.Range("I21").Value = 2646
$ws.Range("J21").Value = 1356
$ws.Range("K21").Value = 2496

# Row 26 - Min 100 purchases: fill in the previously-missing raw counts
$ws.Range("D26").Value = 706
$ws.Range("E26").Value = 256
$ws.Range("F26").Value = 653
$ws.Range("G26").Value = 309
$ws.Range("H26").Value = 626
$ws.Range("I26").Value = 229
$ws.Range("J26").Value = 86
$ws.Range("K26").Value = 6662

[void]$ws.Range("D27").Select()
